# Add example values to the "MIMS" template sheet (row 2 of the annotation table).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MIMS")

$ws.Range("B2").Value = "eukaryotic soil metagenome"
$ws.Range("C2").Value = "NCBITaxon"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/NCBITaxon_2977533"
# "2024-07-10" looks like a date, so force text entry and strip the
# resulting number-format style back off so the cell stays a plain string.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2024-07-10"
$ws.Range("E2").ClearFormats()
$ws.Range("H2").Value = "forest biome"
$ws.Range("I2").Value = "ENVO"
$ws.Range("J2").Value = "https://bioregistry.io/ENVO:01000174"
$ws.Range("K2").Value = "coniferous forest biome"
$ws.Range("L2").Value = "ENVO"
$ws.Range("M2").Value = "https://bioregistry.io/ENVO:01000196"
$ws.Range("N2").Value = "soil"
$ws.Range("O2").Value = "ENVO"
$ws.Range("P2").Value = "https://bioregistry.io/ENVO:00001998"
$ws.Range("Q2").Value = "Germany"
$ws.Range("R2").Value = "NCIT"
$ws.Range("S2").Value = "https://bioregistry.io/NCIT:C16636"
$ws.Range("T2").Value = "+50.55° / +6.21°"
$ws.Range("W2").Value = "Illumina MiSeq"
$ws.Range("X2").Value = "EFO"
$ws.Range("Y2").Value = "https://bioregistry.io/EFO:0004205"
